$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.195.04"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.890.55"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.14"
$ws.Range("E5").Value = "  +8.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.42"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.612"
$ws.Range("E7").Value = "  -1.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -3.07%  "
$ws.Range("E10").Value = "  -5.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000332"
$ws.Range("E11").Value = "  -5.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.90"
$ws.Range("E12").Value = "  -2.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.513.51"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.903.10"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.98"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.21"
$ws.Range("E17").Value = "  +6.64%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.134"
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.82"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.146.61"
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "423.15"
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("E22").Value = "  -5.22%  "
$ws.Range("E24").Value = "  -2.07%  "
$ws.Range("E25").Value = "  +8.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.38"
$ws.Range("E26").Value = "  -7.63%  "
$ws.Range("E27").Value = "  -3.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.26"
$ws.Range("E28").Value = "  -2.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "691.66"
$ws.Range("E29").Value = "  -3.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.18"
$ws.Range("E30").Value = "  -1.46%  "
$ws.Range("E31").Value = "  -3.32%  "
$ws.Range("E32").Value = "  -2.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "68.07"
$ws.Range("E33").Value = "  +10.46%  "
$ws.Range("E34").Value = "  +6.81%  "
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0855"
$ws.Range("E35").Value = "  -4.31%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.89"
$ws.Range("E36").Value = "  -2.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "39.91"
$ws.Range("E37").Value = "  -2.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.147"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.29"
$ws.Range("E41").Value = "  +7.01%  "
$ws.Range("E42").Value = "  +7.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0480"
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.77"
$ws.Range("E44").Value = "  -6.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.39"
$ws.Range("E45").Value = "  +1.75%  "
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("E47").Value = "  +6.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.755.51"
$ws.Range("E48").Value = "  +14.85%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000270"
$ws.Range("E49").Value = "  +8.14%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.42"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.51"
$ws.Range("E51").Value = "  +6.12%  "
